$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Change 1: the "THU Nov 02" / " 10:52:55 PDT 2017" timestamp was split
# across two runs; collapse it into a single run with identical text.
# --------------------------------------------------------------------
$d.Content.Find.Execute("THU Nov 02 10:52:55 PDT 2017", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "THU Nov 02 10:52:55 PDT 2017", 2) | Out-Null

# --------------------------------------------------------------------
# Change 2: append a brand-new purchase record (SAT Nov 04 / BEET /
# MAMATHA CHICK IN) right after the "Amount balance ... - 52137.0"
# paragraph, ahead of the trailing blank paragraphs.
# --------------------------------------------------------------------

$tab = [char]9

# Anchor: the "Amount balance" paragraph that currently ends with "- 52137.0"
$anchorPara = $d.Paragraphs.Item(29)

# 1) New blank bold paragraph directly under the anchor (inherits Bold from it).
$anchorPara.Range.InsertParagraphAfter()

# 2) Nine more blank paragraphs, inserted ahead of the (still) original first
#    trailing blank paragraph so they land, in order, right after the one
#    created above and before the pre-existing blank paragraphs.
for ($i = 0; $i -lt 9; $i++) {
    $d.Paragraphs.Item(31).Range.InsertParagraphBefore()
}

# Paragraph 30: blank, bold (formatting already correct by inheritance).
$d.Paragraphs.Item(30).Range.Font.Bold = $true

# Paragraph 31: "SAT Nov 04" + " 10:18:20 PST 2017" kept as two runs, as in
# the source document's own date-time paragraphs.
$p31 = $d.Paragraphs.Item(31)
$p31.Range.Font.Bold = $false
$start = $p31.Range.Start
$d.Range($start, $start).InsertAfter("SAT Nov 04")
$placeholder = $d.Range($start + 10, $start + 10)
$placeholder.InsertAfter("X")
$boldPlaceholder = $d.Range($start + 10, $start + 11)
$boldPlaceholder.Font.Bold = $true
$boldPlaceholder.Text = " 10:18:20 PST 2017"
$tail = $d.Range($start + 10, $start + 10 + " 10:18:20 PST 2017".Length)
$tail.Font.Bold = $false

# Paragraph 32: "Person Name" ... "- TNP"
$p32 = $d.Paragraphs.Item(32)
$p32.Range.Font.Bold = $false
$p32.Range.Text = "Person Name" + $tab + $tab + $tab + $tab + "- TNP"

# Paragraph 33: dashed separator line
$p33 = $d.Paragraphs.Item(33)
$p33.Range.Font.Bold = $false
$p33.Range.Text = "---------------------------------------------------------------"

# Paragraph 34: "Item Name" ... "- BEET"
$p34 = $d.Paragraphs.Item(34)
$p34.Range.Font.Bold = $false
$p34.Range.Text = "Item Name" + $tab + $tab + $tab + $tab + "- BEET"

# Paragraph 35: "Amount Received" ... "- 1000", in red.
$p35 = $d.Paragraphs.Item(35)
$p35.Range.Font.Bold = $false
$p35.Range.Text = "Amount Received" + $tab + $tab + $tab + "- 1000"
$p35.Range.Font.Color = 255

# Paragraph 36: "Amount balance" ... "- 51137.0", bold.
$p36 = $d.Paragraphs.Item(36)
$p36.Range.Text = "Amount balance" + $tab + $tab + $tab + "- 51137.0"
$p36.Range.Font.Bold = $true

# Paragraph 37: "Amount Received mode" ... "- CASH"
$p37 = $d.Paragraphs.Item(37)
$p37.Range.Font.Bold = $false
$p37.Range.Text = "Amount Received mode" + $tab + $tab + "- CASH"

# Paragraph 38: blank, regular.
$d.Paragraphs.Item(38).Range.Font.Bold = $false

# Paragraph 39: blank, bold.
$d.Paragraphs.Item(39).Range.Font.Bold = $true

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
